# 1/16 diary added --Weihuan Fu
# Applies the content + formatting changes described by the commit diff:
#  - Row 11 (1/12 entry): Reflection/Mood text updated to the Fork-GUI reflection
#  - Row 12 becomes a real data row for the 1/16/2020 class session (was a
#    placeholder/instructions row before)
#  - Row 13 keeps its "Etc." marker
#  - Header row + data rows 9-13 get a thin box border and centered
#    (horizontal+vertical) wrapped text
#  - Columns A-E get narrower/custom widths so the new long-form text fits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Content updates
# (write order mirrors the authoring order so the shared-string table
#  comes out in the same sequence as the saved workbook)
# ---------------------------------------------------------------------

# Row 12 (new 1/16/2020 entry) -- replaces the old placeholder/instructions row
$ws.Range("A12").Value = 43846
$ws.Range("B12").Value = "class time"
$ws.Range("C12").Value = "Xiaolue Peng"
$ws.Range("D12").Value = "understand and practice how to precisely locate issues that needs to be fixed in a system with a bunch of codes and classes"

# Row 11 (1/12/2020 entry) -- Reflection / Goal / Mood text rewritten
$ws.Range("E11").Value = "Ended up having to use a Git GUI App called Fork for cloning, committing, pushing and create merge request. Much easier than use git command line functions to manipulate version control locally"
$ws.Range("F11").Value = "understand the importance of version control especially when everyone if working on different parts of one system. Got familiar with the basic functions of git GUI but more terms need to be learned in the future such as how to deal with merge conflict"
$ws.Range("G11").Value = "excited to learn something that is widely used in the industry."

# Back to row 12 for the remaining reflection columns
$ws.Range("E12").Value = "learned and practiced multiple theories of how to locate issues. Leaned about some common ways of naming files such as ""sprite"". Learn the way of approaching and finally locate the code we want to change"
$ws.Range("F12").Value = "It is hard to understand every line of code in a big system but it is almost unnecessary to be able to do so (because everyone usually has their own modules), however it is very useful to understand how the majority functions works in terms of the software behavior. And being able to quickly do so, locate problem and fix it will be very efficient"
$ws.Range("G12").Value = "Not really familiar with reading codes that are writtened by others but excited to learn a lot from reading other people code."

# Row 13 keeps "Etc." in A13 (unchanged content; B13:G13 stay blank)

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.666666666666668
$ws.Columns.Item(3).ColumnWidth = 25.166666666666668
$ws.Columns.Item(4).ColumnWidth = 41.998697916666664
$ws.Columns.Item(5).ColumnWidth = 40.666666666666664
# Columns F:G keep their original 34.6640625 width

# ---------------------------------------------------------------------
# Row heights (the new long-form text needs extra vertical room)
# ---------------------------------------------------------------------
$ws.Rows.Item(11).RowHeight = 119
$ws.Rows.Item(12).RowHeight = 170

# ---------------------------------------------------------------------
# Formatting: thin box border + centered (H+V) wrapped text across the
# header row and all the data/template rows (9-13)
# ---------------------------------------------------------------------
$tableRange = $ws.Range("A9:G13")
$tableRange.WrapText = $true
$tableRange.HorizontalAlignment = -4108   # xlCenter
$tableRange.VerticalAlignment = -4108     # xlCenter
$tableRange.Borders.LineStyle = 1         # xlContinuous
$tableRange.Borders.Weight = 2            # xlThin

# Active-cell selection, matching the saved workbook state
$ws.Range("F12").Select()
